$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new values look numeric ("1234", "123") but must be stored as text,
# same as the existing rows (e.g. the "123" password/user entries already
# on the sheet). Mark those cells as Text before writing them so they are
# kept as strings instead of being coerced into numbers.
$ws.Range("A5").NumberFormat = "@"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B6").NumberFormat = "@"

# Row 5: new user "1234"
$ws.Range("A5").Value = "1234"
$ws.Range("B5").Value = "123"
$ws.Range("C5").Value = "Cliente"

# Row 6: new user "micha"
$ws.Range("A6").Value = "micha"
$ws.Range("B6").Value = "123"
$ws.Range("C6").Value = "Cliente"
